$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 162; this shifts the existing rows 162:201 down to 163:202
$ws.Rows("162:162").Insert()

# Populate the new row 162 with the new data record
$ws.Range("A162").Value = 11
$ws.Range("B162").Value = "Vega Monumental Concepción"
$ws.Range("C162").Value = "Bíobío"
$ws.Range("D162").Value = 45244
$ws.Range("E162").Value = 8
$ws.Range("F162").Value = 100112001
$ws.Range("G162").Value = "Berenjena"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 150
$ws.Range("K162").Value = 10000
$ws.Range("L162").Value = 10000
$ws.Range("M162").Value = 10000
$ws.Range("N162").Value = "$/caja 50 unidades"
$ws.Range("O162").Value = "Región de Arica y Parinacota"
$ws.Range("P162").Value = 200
$ws.Range("Q162").Value = 50
$ws.Range("R162").Value = "Hortaliza"
